$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.326.14"
$ws.Range("E2").Value = "  +1.03%  "

$ws.Range("D3").Value = "'1.665.86"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("E4").Value = "  +0.84%  "

$ws.Range("D5").Value = "'219.22"
$ws.Range("E5").Value = "  +0.79%  "

$ws.Range("D6").Value = "'0.5349"
$ws.Range("E6").Value = "  +1.70%  "

$ws.Range("E7").Value = "  +0.80%  "

$ws.Range("D8").Value = "'0.2667"
$ws.Range("E8").Value = "  +2.68%  "

$ws.Range("D9").Value = "'0.06403"
$ws.Range("E9").Value = "  +1.19%  "

$ws.Range("D10").Value = "'20.75"
$ws.Range("E10").Value = "  +1.94%  "

$ws.Range("D11").Value = "'0.07855"
$ws.Range("E11").Value = "  +0.61%  "

$ws.Range("D12").Value = "'4.569"
$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").Value = "'1.666.64"
$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("D14").Value = "'1.893.48"
$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("D15").Value = "'0.5537"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").Value = "'0.0₅8225"
$ws.Range("E16").Value = "  +0.11%  "

$ws.Range("E18").Value = "  +0.83%  "

$ws.Range("D19").Value = "'4.692"
$ws.Range("E19").Value = "  +2.62%  "

$ws.Range("D20").Value = "'193.70"
$ws.Range("E20").Value = "  +1.27%  "

$ws.Range("E21").Value = "  +2.18%  "

$ws.Range("D22").Value = "'6.047"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").Value = "'1.012"

$ws.Range("D24").Value = "'146.45"
$ws.Range("E24").Value = "  +2.87%  "

$ws.Range("E25").Value = "  -0.32%  "

$ws.Range("D26").Value = "'7.214"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").Value = "'16.13"

$ws.Range("D28").Value = "'1.504"
$ws.Range("E28").Value = "  +5.27%  "

$ws.Range("D29").Value = "'0.05838"
$ws.Range("E29").Value = "  +0.46%  "

$ws.Range("D30").Value = "'1.283"
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").Value = "'3.641"
$ws.Range("E31").Value = "  +2.80%  "

$ws.Range("D32").Value = "'3.281"
$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("E33").Value = "  +1.96%  "

$ws.Range("D34").Value = "'0.9701"

$ws.Range("D35").Value = "'2.829"
$ws.Range("E35").Value = "  +1.90%  "

$ws.Range("D36").Value = "'2.420"

$ws.Range("D37").Value = "'0.5827"
$ws.Range("E37").Value = "  +1.81%  "

$ws.Range("D38").Value = "'0.01606"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").Value = "'0.8713"
$ws.Range("E39").Value = "  +3.28%  "

$ws.Range("D40").Value = "'5.864"
$ws.Range("E40").Value = "  +2.13%  "

$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "'105.36"
$ws.Range("E41").Value = "  +1.83%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'1.053.49"
$ws.Range("E42").Value = "  +2.62%  "

$ws.Range("D43").Value = "'1.011"
$ws.Range("E43").Value = "  +0.83%  "

$ws.Range("D44").Value = "'1.804.53"
$ws.Range("E44").Value = "  +0.57%  "

$ws.Range("D45").Value = "'57.88"
$ws.Range("E45").Value = "  +1.56%  "

$ws.Range("E46").Value = "  +1.34%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.4388"
$ws.Range("E47").Value = "  +1.71%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.032"
$ws.Range("E48").Value = "  +2.23%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05169"
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.0₈100"
$ws.Range("E50").Value = "  -10.56%  "

$ws.Range("D51").Value = "'1.415"
$ws.Range("E51").Value = "  -3.44%  "
